# Update scheduled market-price data on the Jenova profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 3512925.2
$ws.Range("I86").Value = 3698.6
$ws.Range("K86").Value = 3698.6
$ws.Range("M86").Value = -2575.6

# Row 88
$ws.Range("H88").Value = 1889.4445
$ws.Range("J88").Value = 2157.7144
$ws.Range("L88").Value = 2157.7144
$ws.Range("N88").Value = -2969.7144

# Row 89
$ws.Range("H89").Value = 3512925.2
$ws.Range("I89").Value = 3698.6
$ws.Range("K89").Value = 18493
$ws.Range("M89").Value = -12877

# Row 91
$ws.Range("H91").Value = 1889.4445
$ws.Range("J91").Value = 2157.7144
$ws.Range("L91").Value = 2157.7144
$ws.Range("N91").Value = -4965.7144

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7399.647
$ws.Range("I32").Value = 7488
$ws.Range("J32").Value = 4484
$ws.Range("K32").Value = 7488
$ws.Range("L32").Value = 4484
$ws.Range("M32").Value = -7201

# Row 61
$ws.Range("H61").Value = 2286.946
$ws.Range("I61").Value = 2063.7188
$ws.Range("K61").Value = 2063.7188
$ws.Range("M61").Value = -1851.7188

# Row 132
$ws.Range("H132").Value = 2023.6232
$ws.Range("I132").Value = 2028
$ws.Range("J132").Value = 1990.25
$ws.Range("K132").Value = 6084
$ws.Range("L132").Value = 5970.75
$ws.Range("M132").Value = -3554
$ws.Range("N132").Value = -11030.75

# Row 136
$ws.Range("H136").Value = 2286.946
$ws.Range("I136").Value = 2063.7188
$ws.Range("K136").Value = 6191.1564
$ws.Range("M136").Value = -3641.1564

$ws = $wb.Worksheets.Item("BSM")
# Row 53
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("N53").ClearContents()

# Row 134
$ws.Range("H134").Value = 28259.121
$ws.Range("I134").Value = 4313.3125
$ws.Range("J134").Value = 113399.78
$ws.Range("K134").Value = 12939.9375
$ws.Range("L134").Value = 340199.34
$ws.Range("M134").Value = -10404.9375
$ws.Range("N134").Value = -345269.34

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 518.375
$ws.Range("I7").Value = 492.7
$ws.Range("K7").Value = 492.7
$ws.Range("M7").Value = -379.7

# Row 31
$ws.Range("H31").Value = 31489.945
$ws.Range("I31").Value = 2595.8462
$ws.Range("K31").Value = 2595.8462
$ws.Range("M31").Value = -2300.8462

# Row 34
$ws.Range("H34").Value = 31489.945
$ws.Range("I34").Value = 2595.8462
$ws.Range("K34").Value = 2595.8462
$ws.Range("M34").Value = -2393.8462

$ws = $wb.Worksheets.Item("CUL")
# Row 69
$ws.Range("H69").Value = 468.75
$ws.Range("I69").Value = 733.3333
$ws.Range("K69").Value = 2199.9999
$ws.Range("M69").Value = -1388.9999

# Row 72
$ws.Range("H72").Value = 468.75
$ws.Range("I72").Value = 733.3333
$ws.Range("K72").Value = 6599.9997
$ws.Range("M72").Value = -2543.9997

# Row 115
$ws.Range("H115").Value = 2904.1667
$ws.Range("I115").Value = 2748.5
$ws.Range("J115").Value = 3215.5
$ws.Range("K115").Value = 8245.5
$ws.Range("L115").Value = 9646.5
$ws.Range("M115").Value = -7070.5
$ws.Range("N115").Value = -11996.5

# Row 131
$ws.Range("H131").Value = 3815.3416
$ws.Range("I131").Value = 1272.7142
$ws.Range("J131").Value = 4338.8237
$ws.Range("K131").Value = 3818.1426
$ws.Range("L131").Value = 13016.4711
$ws.Range("M131").Value = 1221.8574
$ws.Range("N131").Value = -23096.4711

$ws = $wb.Worksheets.Item("GSM")
# Row 99
$ws.Range("H99").Value = 4735
$ws.Range("J99").Value = 4999
$ws.Range("L99").Value = 4999

# Row 122
$ws.Range("H122").Value = 6572.4287
$ws.Range("I122").Value = 7003.5
$ws.Range("K122").Value = 21010.5
$ws.Range("M122").Value = -18560.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2010740.8
$ws.Range("I7").Value = 30004
$ws.Range("J7").Value = 2505925
$ws.Range("K7").Value = 30004
$ws.Range("L7").Value = 2505925
$ws.Range("M7").Value = -29892
$ws.Range("N7").Value = -2506149

# Row 40
$ws.Range("H40").Value = 338668
$ws.Range("I40").Value = 1000004
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 1000004
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -999868
$ws.Range("N40").Value = -8272

# Row 46
$ws.Range("H46").Value = 3958.96
$ws.Range("J46").Value = 5745.8335
$ws.Range("L46").Value = 5745.8335
$ws.Range("N46").Value = -6121.8335

# Row 122
$ws.Range("H122").Value = 1113403.2
$ws.Range("I122").Value = 627422.4
$ws.Range("J122").Value = 5001250
$ws.Range("K122").Value = 1882267.2
$ws.Range("L122").Value = 15003750
$ws.Range("M122").Value = -1879817.2
$ws.Range("N122").Value = -15008650

# Row 126
$ws.Range("H126").Value = 2010740.8
$ws.Range("I126").Value = 30004
$ws.Range("J126").Value = 2505925
$ws.Range("K126").Value = 90012
$ws.Range("L126").Value = 7517775
$ws.Range("M126").Value = -87542
$ws.Range("N126").Value = -7522715

# Row 132
$ws.Range("H132").Value = 6378.2
$ws.Range("I132").Value = 5871.316
$ws.Range("J132").Value = 7983.3335
$ws.Range("K132").Value = 17613.948
$ws.Range("L132").Value = 23950.0005
$ws.Range("M132").Value = -15083.948
$ws.Range("N132").Value = -29010.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 100000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

# Row 67
$ws.Range("H67").Value = 100000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

# Row 107
$ws.Range("H107").Value = 893.5789
$ws.Range("I107").Value = 1079.4615
$ws.Range("J107").Value = 490.83334
$ws.Range("K107").Value = 3238.3845
$ws.Range("L107").Value = 1472.50002
$ws.Range("M107").Value = -1318.3845
$ws.Range("N107").Value = -5312.500019999999

# Row 118
$ws.Range("H118").Value = 49990
$ws.Range("J118").Value = 49990
$ws.Range("L118").Value = 49990
$ws.Range("N118").Value = -53304

# Row 122
$ws.Range("H122").Value = 38463404
$ws.Range("I122").Value = 50001310
$ws.Range("K122").Value = 150003930
$ws.Range("M122").Value = -150001480

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("M126").ClearContents()
